$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.243.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.51%  "
$ws.Range("D3").Value = "'1.905.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.40%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'307.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.66%  "
$ws.Range("D6").Value = "'0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").Value = "'0.5239"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.19%  "
$ws.Range("D8").Value = "'0.3779"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.51%  "
$ws.Range("D9").Value = "'0.07254"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.37%  "
$ws.Range("D10").Value = "'21.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.65%  "
$ws.Range("D11").Value = "'0.8996"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.27%  "
$ws.Range("D12").Value = "'0.08265"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +10.48%  "
$ws.Range("D13").Value = "'1.909.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.55%  "
$ws.Range("D14").Value = "'95.35"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.37%  "
$ws.Range("D15").Value = "'5.279"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.34%  "
$ws.Range("D16").Value = "'1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").Value = "'0.000008603"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.46%  "
$ws.Range("D18").Value = "'14.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.47%  "
$ws.Range("D19").Value = "'0.9999"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("D20").Value = "'27.280.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.56%  "
$ws.Range("D21").Value = "'5.064"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.55%  "
$ws.Range("D22").Value = "'2.151.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.43%  "
$ws.Range("D23").Value = "'10.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.23%  "
$ws.Range("D24").Value = "'6.454"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.66%  "
$ws.Range("D25").Value = "'2.300"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.13%  "
$ws.Range("D26").Value = "'145.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.69%  "
$ws.Range("D27").Value = "'1.744"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.43%  "
$ws.Range("D28").Value = "'18.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.86%  "
$ws.Range("D29").Value = "'114.83"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.20%  "
$ws.Range("D30").Value = "'4.977"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.92%  "
$ws.Range("D31").Value = "'4.813"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.08%  "
$ws.Range("D32").Value = "'0.09199"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.91%  "
$ws.Range("D33").Value = "'0.8037"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.64%  "
$ws.Range("D34").Value = "'0.05083"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.12%  "
$ws.Range("D35").Value = "'1.239"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.85%  "
$ws.Range("D36").Value = "'2.942"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("D37").Value = "'3.349"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.50%  "
$ws.Range("D38").Value = "'2.572"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.11%  "
$ws.Range("D39").Value = "'0.5732"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.60%  "
$ws.Range("D40").Value = "'0.01979"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("D41").Value = "'1.074"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.31%  "
$ws.Range("D42").Value = "'9.064"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.91%  "
$ws.Range("D43").Value = "'6.633"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.06%  "
$ws.Range("D44").Value = "'118.60"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.37%  "
$ws.Range("D45").Value = "'0.1518"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.48%  "
$ws.Range("D46").Value = "'0.4843"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.35%  "
$ws.Range("D47").Value = "'10.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.74%  "
$ws.Range("D48").Value = "'0.9998"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("D49").Value = "'1.612"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.08%  "
$ws.Range("D50").Value = "'37.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.46%  "
$ws.Range("D51").Value = "'63.69"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.46%  "

# Rows 47 and 48 swap Coin name and Link (data reordering upstream)
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
